# Added check for correct unit type while importing fund units
#
# The sample "fund_units" import sheet contained an invalid Unit Type
# value ("Series C") in the Unit Type column (D) that the import
# validation doesn't accept. Correct the sample data so every row uses
# one of the valid unit types (Series A / Series B), and give the
# corrected column its own ("Normal 3") cell style so it's visibly
# distinct/reviewed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Unit Type *" column (D) ------------------------------
# Rows 2-6  (Call 1 block)          Rows 7-11 (Distribution 1 block)
#   D2 Series C -> Series A           D7  Series C -> Series A
#   D3 Series C -> Series A           D8  Series C -> Series A
#   D4 Series A -> Series B           D9  Series A -> Series B
#   D5 Series A -> Series B           D10 Series A -> Series B
#   D6 Series B -> Series B           D11 Series B -> Series B
$ws.Range("D2").Value = "Series A"
$ws.Range("D3").Value = "Series A"
$ws.Range("D4").Value = "Series B"
$ws.Range("D5").Value = "Series B"
$ws.Range("D6").Value = "Series B"

$ws.Range("D7").Value = "Series A"
$ws.Range("D8").Value = "Series A"
$ws.Range("D9").Value = "Series B"
$ws.Range("D10").Value = "Series B"
$ws.Range("D11").Value = "Series B"

# --- Apply the "Normal 3" style to the corrected cells --------------
$ws.Range("D2:D6").Style = "Normal 3"
$ws.Range("D7:D11").Style = "Normal 3"

# --- Leave the selection on the corrected range ----------------------
$null = $ws.Range("D7:D11").Select()
